$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: re_rank score updated
$ws.Range("F2").Value = 8.43712372573634

# Row 3: re_rank score updated
$ws.Range("F3").Value = 8.085171174491482

# Row 4 and Row 5 swap identities (prolificid/name/gender) and B (prolificid-order) swaps
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("D4").Value = "Alfredo"
$ws.Range("E4").Value = "male"
$ws.Range("F4").Value = 7.483015296297952

$ws.Range("B5").Value = 8
$ws.Range("C5").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("D5").Value = "Valeria"
$ws.Range("E5").Value = "female"
$ws.Range("F5").Value = 7.18543091573438

# Row 6-8: re_rank score updated
$ws.Range("F6").Value = 6.391489871176138
$ws.Range("F7").Value = 6.387011644639443
$ws.Range("F8").Value = 5.012348235563821

# Row 9 and Row 10 swap identities (prolificid/name) and B swaps
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "5e706891c396cc64388ef760"
$ws.Range("D9").Value = "Maria"
$ws.Range("F9").Value = 3.142757349846526

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "5e0adc8f4cac6834756db412"
$ws.Range("D10").Value = "Mary"
$ws.Range("F10").Value = 3.037402404511541

# Row 11-13: re_rank score updated
$ws.Range("F11").Value = 2.344227293246886
$ws.Range("F12").Value = 1.337899327810482
$ws.Range("F13").Value = 0.002543171126171584
